# Update dashboards - 2026-01-08
# Applies the periodic data refresh described by the commit: new release
# dates (highlighted yellow for newly-updated series) and the corresponding
# shifted trailing-observation values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Cells whose release date just moved onto a "new data" cadence get the
#    yellow highlighted date style (same look already used by N29/N30/N47-50).
#    Copy number-format + fill from an existing "highlighted date" cell so we
#    reuse the workbook's existing style definition instead of inventing one.
# ---------------------------------------------------------------------------
$styleTemplate = $ws.Range("N29")
$dateCellsToHighlight = @("N10", "N11", "N12", "N13", "N14", "C35")

foreach ($cellRef in $dateCellsToHighlight) {
    $styleTemplate.Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Write the refreshed values (release dates + the five trailing
#    observations / growth figures that rolled forward one period).
# ---------------------------------------------------------------------------
$values = @{
    "N10" = 45962; "Q10" = 4.3;   "R10" = 4.5;   "S10" = 4.6;   "U10" = 4.3
    "N11" = 45962; "S11" = 3.4;   "T11" = 3.2
    "N12" = 45962; "R12" = 3.2;   "S12" = 3.3;   "T12" = 3.2
    "N13" = 46020; "Q13" = 208000; "R13" = 200000; "S13" = 215000; "T13" = 224000; "U13" = 237000
    "N14" = 46013; "Q14" = 1914000; "R14" = 1858000; "S14" = 1914000; "T14" = 1885000; "U14" = 1830000

    "F28" = -0.02191649132412532; "G28" = 0.006436255758670795

    "F29" = 0.04772459132664544;  "G29" = 0.07412067603746038
    "N29" = 46029; "R29" = 2.24;  "S29" = 2.23;  "T29" = 2.22;  "U29" = 2.24

    "F30" = -0.01529652492391287; "G30" = 0.001174064535676367
    "N30" = 46029; "R30" = 2.27;  "S30" = 2.26;  "U30" = 2.25

    "F31" = 0.04749518938811943;  "G31" = 0.06502168244015354

    "C35" = 45839; "F35" = 4.9; "G35" = 4.1; "H35" = -2.1; "I35" = 0.9; "J35" = 3.1

    "N47" = 46028

    "N48" = 46028; "Q48" = 3.47; "R48" = 3.46; "T48" = 3.47

    "N49" = 46028; "Q49" = 3.72; "R49" = 3.71; "S49" = 3.74; "T49" = 3.73; "U49" = 3.68

    "N50" = 46028; "Q50" = 4.18; "R50" = 4.17; "S50" = 4.19; "T50" = 4.18; "U50" = 4.14

    "N52" = 46028; "R52" = 5.92; "S52" = 5.93; "T52" = 5.9;  "U52" = 5.89
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
